# Applies the per-coin price/volume refresh from the Oct 28 2024
# GitHub Actions data pull, including the Aptos/SuiNetwork row swap
# and the Cronos -> Optimism row-51 replacement.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # The source data stores every cell (even numeric-looking ones
    # like "601.45") as literal text. Flip to the Text number
    # format before writing so Excel does not auto-convert the
    # string into a real number, then clear the format override
    # again so the cell keeps the workbooks default (General) style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "69.734.46"
Set-TextValue "E2" "  +2.47%  "

Set-TextValue "D3" "2.540.75"
Set-TextValue "E3" "  +1.01%  "

Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.14%  "

Set-TextValue "D5" "601.45"
Set-TextValue "E5" "  +1.99%  "

Set-TextValue "D6" "177.52"
Set-TextValue "E6" "  -0.22%  "

Set-TextValue "E7" "  -0.03%  "

Set-TextValue "D8" "0.519"
Set-TextValue "E8" "  +0.63%  "

Set-TextValue "D9" "2.540.46"
Set-TextValue "E9" "  +1.33%  "

Set-TextValue "E10" "  +12.45%  "

Set-TextValue "E11" "  -0.23%  "

Set-TextValue "E12" "  +1.14%  "

Set-TextValue "D13" "5.03"
Set-TextValue "E13" "  +1.48%  "

Set-TextValue "D14" "0.0000182"
Set-TextValue "E14" "  +5.68%  "

Set-TextValue "D15" "2.982.42"
Set-TextValue "E15" "  +1.92%  "

Set-TextValue "D16" "26.13"
Set-TextValue "E16" "  +1.37%  "

Set-TextValue "D17" "69.656.81"
Set-TextValue "E17" "  +2.62%  "

Set-TextValue "D18" "2.534.72"
Set-TextValue "E18" "  +1.93%  "

Set-TextValue "D19" "7.70"
Set-TextValue "E19" "  +1.97%  "

Set-TextValue "D20" "365.02"
Set-TextValue "E20" "  +3.21%  "

Set-TextValue "D21" "11.11"
Set-TextValue "E21" "  +0.68%  "

Set-TextValue "D22" "4.10"

Set-TextValue "E23" "  -0.17%  "

Set-TextValue "E24" "  -0.38%  "

Set-TextValue "D25" "4.27"
Set-TextValue "E25" "  -1.73%  "

Set-TextValue "B26" "Aptos"
Set-TextValue "C26" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D26" "9.20"
Set-TextValue "E26" "  +0.48%  "

Set-TextValue "B27" "SuiNetwork"
Set-TextValue "C27" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D27" "1.72"
Set-TextValue "E27" "  -2.84%  "

Set-TextValue "D28" "2.667.01"
Set-TextValue "E28" "  +2.89%  "

Set-TextValue "D29" "0.997"
Set-TextValue "E29" "  +0.43%  "

Set-TextValue "E30" "  +0.44%  "

Set-TextValue "D31" "513.95"
Set-TextValue "E31" "  +0.93%  "

Set-TextValue "D32" "7.78"
Set-TextValue "E32" "  -1.32%  "

Set-TextValue "D33" "1.26"
Set-TextValue "E33" "  -0.98%  "

Set-TextValue "D34" "1.80"
Set-TextValue "E34" "  +1.11%  "

Set-TextValue "E35" "  +0.17%  "

Set-TextValue "D36" "0.120"
Set-TextValue "E36" "  -1.73%  "

Set-TextValue "D37" "161.53"
Set-TextValue "E37" "  -1.92%  "

Set-TextValue "D38" "18.82"
Set-TextValue "E38" "  +2.03%  "

Set-TextValue "D39" "18.91"
Set-TextValue "E39" "  +1.36%  "

Set-TextValue "E40" "  -0.95%  "

Set-TextValue "E41" "  -0.01%  "

Set-TextValue "D42" "1.75"
Set-TextValue "E42" "  +0.13%  "

Set-TextValue "D43" "4.88"
Set-TextValue "E43" "  -0.61%  "

Set-TextValue "D44" "0.323"
Set-TextValue "E44" "  -2.32%  "

Set-TextValue "D45" "2.45"
Set-TextValue "E45" "  -1.69%  "

Set-TextValue "D46" "38.76"
Set-TextValue "E46" "  -0.53%  "

Set-TextValue "D47" "151.38"
Set-TextValue "E47" "  +3.86%  "

Set-TextValue "D48" "3.62"
Set-TextValue "E48" "  +2.02%  "

Set-TextValue "D49" "0.521"
Set-TextValue "E49" "  -0.15%  "

Set-TextValue "D50" "0.0₆0255"
Set-TextValue "E50" "  -1.71%  "

Set-TextValue "B51" "Optimism"
Set-TextValue "C51" "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextValue "D51" "1.60"
Set-TextValue "E51" "  +0.08%  "
